# Update the two-digit x two-digit multiplication problems in the table.
# Each content row (1, 5, 10, 15, 20) has 5 problem cells (columns 1-5).
# Cell.Range.Text preserves the existing run/paragraph formatting
# (font, size, justification) already present in each cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="61×39="},
    @{Row=1;  Col=2; Text="21×94="},
    @{Row=1;  Col=3; Text="11×48="},
    @{Row=1;  Col=4; Text="85×62="},
    @{Row=1;  Col=5; Text="58×93="},

    @{Row=5;  Col=1; Text="75×17="},
    @{Row=5;  Col=2; Text="50×27="},
    @{Row=5;  Col=3; Text="65×97="},
    @{Row=5;  Col=4; Text="81×94="},
    @{Row=5;  Col=5; Text="12×72="},

    @{Row=10; Col=1; Text="93×16="},
    @{Row=10; Col=2; Text="81×42="},
    @{Row=10; Col=3; Text="98×33="},
    @{Row=10; Col=4; Text="64×18="},
    @{Row=10; Col=5; Text="27×88="},

    @{Row=15; Col=1; Text="25×81="},
    @{Row=15; Col=2; Text="15×28="},
    @{Row=15; Col=3; Text="42×47="},
    @{Row=15; Col=4; Text="70×89="},
    @{Row=15; Col=5; Text="64×58="},

    @{Row=20; Col=1; Text="98×75="},
    @{Row=20; Col=2; Text="48×45="},
    @{Row=20; Col=3; Text="91×14="},
    @{Row=20; Col=4; Text="18×64="},
    @{Row=20; Col=5; Text="62×89="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
